$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Goodreads" row before current row 21 (Amazon), matching sort by HMOD (column D) descending.
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21,1).Value = "Goodreads"
$ws.Cells.Item(21,2).Value = "5A471B"
$ws.Cells.Item(21,3).Value = 42
$ws.Cells.Item(21,4).Formula = "=MOD((C21+100),360)"
$ws.Cells.Item(21,5).Value = 70
$ws.Cells.Item(21,6).Value = 35

# Insert "Soundcloud" row after RSS (now row 24), before Stumbleupon.
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25,1).Value = "Soundcloud"
$ws.Cells.Item(25,2).Value = "FF6600"
$ws.Cells.Item(25,3).Value = 24
$ws.Cells.Item(25,4).Formula = "=MOD((C25+100),360)"
$ws.Cells.Item(25,5).Value = 100
$ws.Cells.Item(25,6).Value = 100

# Insert "Email" row after 500px (now row 30), before Pinterest.
$ws.Rows.Item(31).Insert()
$ws.Cells.Item(31,1).Value = "Email"
$ws.Cells.Item(31,2).Value = 666666
$ws.Cells.Item(31,3).Value = 0
$ws.Cells.Item(31,4).Formula = "=MOD((C31+100),360)"
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 40

# Update dimension/selection/sort-state references to the new extent (A1:F34).
$ws.Range("A32").Select()

$wbView = $excel.ActiveWindow
